$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: clear the stray Hydrogen/Non-metallic-minerals value (was 1297.543990558612)
$ws.Range("D3").Value = ""

# C4: corrected Methanol/Chemicals figure
$ws.Range("C4").Value = 50.55362508600344

# C5: corrected Ammonia/Chemicals figure
$ws.Range("C5").Value = 3859.94789163173

# Row 7 label changes from "Other" to "Biogas", and gets its own D-column value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 570.9193558457891

# New row 8 ("Other") appended below, re-using row 7's label formatting
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 0

$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
